$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 111.51282
$ws.Range("I33").Value = 101.97143
$ws.Range("K33").Value = 101.97143
$ws.Range("M33").Value = 127.02857

$ws.Range("H113").Value = 2526
$ws.Range("I113").Value = 2170
$ws.Range("J113").Value = 3000.6667
$ws.Range("K113").Value = 2170
$ws.Range("L113").Value = 3000.6667
$ws.Range("M113").Value = 1084
$ws.Range("N113").Value = -9508.6667

$ws.Range("H116").Value = 7145306.5
$ws.Range("I116").Value = 15386562
$ws.Range("J116").Value = 2885.1333
$ws.Range("K116").Value = 15386562
$ws.Range("L116").Value = 2885.1333
$ws.Range("M116").Value = -15383120
$ws.Range("N116").Value = -9769.1333

$ws.Range("H129").Value = 880.58826
$ws.Range("I129").Value = 529
$ws.Range("J129").Value = 955.9286
$ws.Range("K129").Value = 1587
$ws.Range("L129").Value = 2867.7858
$ws.Range("M129").Value = 3413
$ws.Range("N129").Value = -12867.7858

$ws.Range("H132").Value = 2378.0833
$ws.Range("I132").Value = 1740.9032
$ws.Range("K132").Value = 5222.7096
$ws.Range("M132").Value = -2692.7096

$ws.Range("H135").Value = 595.85364
$ws.Range("I135").Value = 564.87177
$ws.Range("K135").Value = 5083.845929999999
$ws.Range("M135").Value = -2548.845929999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14741.659
$ws.Range("I32").Value = 15764.846
$ws.Range("K32").Value = 15764.846
$ws.Range("M32").Value = -15477.846

$ws.Range("H61").Value = 2754.875
$ws.Range("I61").Value = 1485.7084
$ws.Range("J61").Value = 6562.375
$ws.Range("K61").Value = 1485.7084
$ws.Range("L61").Value = 6562.375
$ws.Range("M61").Value = -1273.7084
$ws.Range("N61").Value = -6986.375

$ws.Range("H74").Value = 1051.9231
$ws.Range("I74").Value = 1193.0385
$ws.Range("J74").Value = 769.6923
$ws.Range("K74").Value = 1193.0385
$ws.Range("L74").Value = 769.6923
$ws.Range("M74").Value = -319.0385000000001
$ws.Range("N74").Value = -2517.6923

$ws.Range("H77").Value = 1051.9231
$ws.Range("I77").Value = 1193.0385
$ws.Range("J77").Value = 769.6923
$ws.Range("K77").Value = 5965.192500000001
$ws.Range("L77").Value = 3848.4615
$ws.Range("M77").Value = -1597.192500000001
$ws.Range("N77").Value = -12584.4615

$ws.Range("H80").Value = 33531.5
$ws.Range("I80").Value = 23000
$ws.Range("J80").Value = 35637.8
$ws.Range("K80").Value = 23000
$ws.Range("L80").Value = 35637.8
$ws.Range("M80").Value = -22002
$ws.Range("N80").Value = -37633.8

$ws.Range("H83").Value = 33531.5
$ws.Range("I83").Value = 23000
$ws.Range("J83").Value = 35637.8
$ws.Range("K83").Value = 69000
$ws.Range("L83").Value = 106913.4
$ws.Range("M83").Value = -64008
$ws.Range("N83").Value = -116897.4

$ws.Range("H132").Value = 3166.5151
$ws.Range("I132").Value = 2768.625
$ws.Range("J132").Value = 3541
$ws.Range("K132").Value = 8305.875
$ws.Range("L132").Value = 10623
$ws.Range("M132").Value = -5775.875
$ws.Range("N132").Value = -15683

$ws.Range("H136").Value = 2754.875
$ws.Range("I136").Value = 1485.7084
$ws.Range("J136").Value = 6562.375
$ws.Range("K136").Value = 4457.1252
$ws.Range("L136").Value = 19687.125
$ws.Range("M136").Value = -1907.1252
$ws.Range("N136").Value = -24787.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 25000
$ws.Range("J58").Value = 25000
$ws.Range("L58").Value = 25000
$ws.Range("N58").Value = -25588

$ws.Range("H99").Value = 1740.35
$ws.Range("I99").Value = 1234.5333
$ws.Range("J99").Value = 3257.8
$ws.Range("K99").Value = 1234.5333
$ws.Range("L99").Value = 3257.8
$ws.Range("M99").Value = 263.4666999999999
$ws.Range("N99").Value = -6253.8

$ws.Range("H134").Value = 2945.1
$ws.Range("I134").Value = 2806.9167
$ws.Range("J134").Value = 3497.8333
$ws.Range("K134").Value = 8420.750100000001
$ws.Range("L134").Value = 10493.4999
$ws.Range("M134").Value = -5885.750100000001
$ws.Range("N134").Value = -15563.4999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 562171.9399999999
$ws.Range("I58").Value = 772460.6
$ws.Range("K58").Value = 772460.6
$ws.Range("M58").Value = -772257.6

$ws.Range("H99").Value = 2807.05
$ws.Range("I99").Value = 2689.3125
$ws.Range("J99").Value = 3278
$ws.Range("K99").Value = 2689.3125
$ws.Range("L99").Value = 3278
$ws.Range("M99").Value = -1191.3125
$ws.Range("N99").Value = -6274

$ws.Range("H105").Value = 11944.444
$ws.Range("I105").Value = 13325
$ws.Range("J105").Value = 900
$ws.Range("K105").Value = 13325
$ws.Range("L105").Value = 900
$ws.Range("M105").Value = -11578
$ws.Range("N105").Value = -4394

$ws.Range("H122").Value = 2665.5264
$ws.Range("I122").Value = 2627.682
$ws.Range("J122").Value = 2717.5625
$ws.Range("K122").Value = 7883.045999999999
$ws.Range("L122").Value = 8152.6875
$ws.Range("M122").Value = -5433.045999999999
$ws.Range("N122").Value = -13052.6875

$ws.Range("H126").Value = 2807.05
$ws.Range("I126").Value = 2689.3125
$ws.Range("J126").Value = 3278
$ws.Range("K126").Value = 8067.9375
$ws.Range("L126").Value = 9834
$ws.Range("M126").Value = -5597.9375
$ws.Range("N126").Value = -14774

$ws.Range("H132").Value = 230277.19
$ws.Range("I132").Value = 265898.44
$ws.Range("J132").Value = 3191.75
$ws.Range("K132").Value = 797695.3200000001
$ws.Range("L132").Value = 9575.25
$ws.Range("M132").Value = -795165.3200000001
$ws.Range("N132").Value = -14635.25

$ws.Range("H134").Value = 1234.6842
$ws.Range("I134").Value = 1102.5676
$ws.Range("J134").Value = 1479.1
$ws.Range("K134").Value = 3307.7028
$ws.Range("L134").Value = 4437.299999999999
$ws.Range("M134").Value = -772.7028
$ws.Range("N134").Value = -9507.299999999999

$ws.Range("H136").Value = 562171.9399999999
$ws.Range("I136").Value = 772460.6
$ws.Range("K136").Value = 2317381.8
$ws.Range("M136").Value = -2314831.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 92.5
$ws.Range("I40").Value = 92.5
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 370
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -301
$ws.Range("N40").ClearContents()

$ws.Range("H137").Value = 3285.5715
$ws.Range("I137").Value = 1300
$ws.Range("J137").Value = 3616.5
$ws.Range("K137").Value = 3900
$ws.Range("L137").Value = 10849.5
$ws.Range("M137").Value = 1200
$ws.Range("N137").Value = -21049.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1801.762
$ws.Range("I132").Value = 1267.48
$ws.Range("J132").Value = 2587.4707
$ws.Range("K132").Value = 3802.44
$ws.Range("L132").Value = 7762.4121
$ws.Range("M132").Value = -1272.44
$ws.Range("N132").Value = -12822.4121

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5037.9165
$ws.Range("I132").Value = 4733.857
$ws.Range("J132").Value = 7166.3335
$ws.Range("K132").Value = 14201.571
$ws.Range("L132").Value = 21499.0005
$ws.Range("M132").Value = -11671.571
$ws.Range("N132").Value = -26559.0005

$ws.Range("H136").Value = 1839.9854
$ws.Range("I136").Value = 1405.5088
$ws.Range("J136").Value = 4091.3635
$ws.Range("K136").Value = 4216.526400000001
$ws.Range("L136").Value = 12274.0905
$ws.Range("M136").Value = -1666.526400000001
$ws.Range("N136").Value = -17374.0905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5269.483
$ws.Range("I126").Value = 5869.2
$ws.Range("K126").Value = 17607.6
$ws.Range("M126").Value = -15137.6

$ws.Range("H136").Value = 1379.5927
$ws.Range("I136").Value = 1449.96
$ws.Range("J136").Value = 500
$ws.Range("K136").Value = 4349.88
$ws.Range("L136").Value = 1500
$ws.Range("M136").Value = -1799.88
$ws.Range("N136").Value = -6600
